# daily auto push: 2026-01-31 13:49 UTC
#
# A new observation row is inserted into the "sei3" log sheet for
# 2026/01/31 (Sat) at time 19, pushing every subsequent row down by one
# (old row 731 "2026/12/29" becomes 732, ... old row 772 becomes 773).
# The sheet's used-range grows from A1:D772 to A1:D773.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 731.. down by one, leaving a blank row 731 to fill in.
$ws.Rows.Item(731).Insert()

# Force column A to be stored as literal text (not auto-parsed into a
# date serial) by using Excel's leading-apostrophe text marker, matching
# how every other date cell in this column is stored.
$ws.Range("A731").Value = "'2026/01/31"
$ws.Range("B731").Value = "土"
$ws.Range("C731").Value = 19
$ws.Range("D731").Value = 22
